# Update the "Förändrad" (changed) date column (C) for rows 2-13 from
# serial date 46060 (2026-02-07) to 46061 (2026-02-08).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46060) {
        $cell.Value = 46061
    }
}
